$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading apostrophe forces text storage for values that would otherwise
# be auto-converted to numbers by Excel (matches original inline-string type).

$ws.Range("D2").Value = "54.430.06"
$ws.Range("E2").Value = "  -7.61%  "

$ws.Range("D3").Value = "2.882.61"
$ws.Range("E3").Value = "  -10.65%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'476.68"
$ws.Range("E5").Value = "  -11.69%  "

$ws.Range("D6").Value = "'126.81"
$ws.Range("E6").Value = "  -6.88%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "2.877.63"
$ws.Range("E8").Value = "  -10.81%  "

$ws.Range("D9").Value = "'0.406"
$ws.Range("E9").Value = "  -11.76%  "

$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  -12.23%  "

$ws.Range("D11").Value = "'0.0976"
$ws.Range("E11").Value = "  -14.94%  "

$ws.Range("E12").Value = "  -15.26%  "

$ws.Range("E13").Value = "  -3.94%  "

$ws.Range("D14").Value = "3.363.34"
$ws.Range("E14").Value = "  -11.06%  "

$ws.Range("D15").Value = "'22.83"
$ws.Range("E15").Value = "  -12.20%  "

$ws.Range("D16").Value = "54.362.45"
$ws.Range("E16").Value = "  -7.83%  "

$ws.Range("D17").Value = "2.873.41"
$ws.Range("E17").Value = "  -11.04%  "

$ws.Range("D18").Value = "'0.0000136"
$ws.Range("E18").Value = "  -14.49%  "

$ws.Range("D19").Value = "'5.23"
$ws.Range("E19").Value = "  -11.78%  "

$ws.Range("D20").Value = "'11.64"
$ws.Range("E20").Value = "  -13.06%  "

$ws.Range("D21").Value = "'7.13"
$ws.Range("E21").Value = "  -13.28%  "

$ws.Range("D22").Value = "'309.57"
$ws.Range("E22").Value = "  -14.73%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -14.06%  "

$ws.Range("D25").Value = "'59.71"
$ws.Range("E25").Value = "  -15.41%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -10.09%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").Value = "0.0₃0823"
$ws.Range("E29").Value = "  -15.10%  "

$ws.Range("D30").Value = "'6.30"
$ws.Range("E30").Value = "  -11.82%  "

$ws.Range("E31").Value = "  -5.94%  "

$ws.Range("D32").Value = "'6.21"
$ws.Range("E32").Value = "  -12.49%  "

$ws.Range("D33").Value = "'19.14"
$ws.Range("E33").Value = "  -12.63%  "

$ws.Range("E34").Value = "  -15.90%  "

$ws.Range("E35").Value = "  -13.86%  "

$ws.Range("D36").Value = "'137.10"
$ws.Range("E36").Value = "  -15.04%  "

$ws.Range("D37").Value = "'5.46"
$ws.Range("E37").Value = "  -15.08%  "

$ws.Range("D38").Value = "'1.22"
$ws.Range("E38").Value = "  -15.44%  "

$ws.Range("D39").Value = "'23.08"
$ws.Range("E39").Value = "  -12.43%  "

$ws.Range("D40").Value = "'0.0622"
$ws.Range("E40").Value = "  -12.28%  "

$ws.Range("D41").Value = "2.901.88"
$ws.Range("E41").Value = "  -10.89%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("D43").Value = "'35.58"
$ws.Range("E43").Value = "  -13.52%  "

$ws.Range("D44").Value = "'0.966"
$ws.Range("E44").Value = "  -13.04%  "

$ws.Range("D45").Value = "'0.602"
$ws.Range("E45").Value = "  -16.02%  "

$ws.Range("D46").Value = "'3.44"
$ws.Range("E46").Value = "  -14.82%  "

$ws.Range("E47").Value = "  -11.98%  "

$ws.Range("D48").Value = "2.057.48"
$ws.Range("E48").Value = "  -10.66%  "

$ws.Range("D49").Value = "'5.34"
$ws.Range("E49").Value = "  -15.51%  "

$ws.Range("D50").Value = "'17.96"
$ws.Range("E50").Value = "  -14.04%  "

$ws.Range("E51").Value = "  -11.79%  "
